# Updated cryptos list values (prices + 1h volume %) per the target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells we touch keep their original text storage
# (many values look numeric, e.g. "1.002" or "0.06760", and Excel would
# otherwise silently coerce them to numbers and normalize/truncate them).

$ws.Range("D2").Value = "30.801.05"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").Value = "2.114.02"
$ws.Range("E3").Value = "  +6.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.53"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5338"
$ws.Range("E7").Value = "  +4.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4395"
$ws.Range("E8").Value = "  +6.72%  "
$ws.Range("E9").Value = "  +2.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.14"
$ws.Range("E10").Value = "  +10.50%  "
$ws.Range("E11").Value = "  +4.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.01"
$ws.Range("E12").Value = "  +3.21%  "
$ws.Range("D13").Value = "2.117.28"
$ws.Range("E13").Value = "  +6.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.763"
$ws.Range("E14").Value = "  +4.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.812"
$ws.Range("E15").Value = "  +5.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.05"
$ws.Range("E16").Value = "  +3.34%  "
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001133"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.14"
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.341"
$ws.Range("E22").Value = "  +4.35%  "
$ws.Range("D23").Value = "30.872.76"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.33"
$ws.Range("E24").Value = "  +7.05%  "
$ws.Range("D25").Value = "2.365.80"
$ws.Range("E25").Value = "  +6.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.286"
$ws.Range("E26").Value = "  +3.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.77"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.598"
$ws.Range("E28").Value = "  +9.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.47"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.34"
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.182"
$ws.Range("E31").Value = "  +3.98%  "
$ws.Range("E32").Value = "  +2.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.253"
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.012"
$ws.Range("E34").Value = "  +5.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.559"
$ws.Range("E35").Value = "  +18.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02603"
$ws.Range("E36").Value = "  +5.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.559"
$ws.Range("E37").Value = "  +3.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "12.88"
$ws.Range("E38").Value = "  +9.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06760"
$ws.Range("E39").Value = "  +4.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.477"
$ws.Range("E40").Value = "  +6.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2287"
$ws.Range("E41").Value = "  +4.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6861"
$ws.Range("E42").Value = "  +4.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.249"
$ws.Range("E43").Value = "  +2.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9995"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.228"
$ws.Range("E47").Value = "  +1.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.666"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.278"
$ws.Range("E49").Value = "  +4.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.19"
$ws.Range("E50").Value = "  +4.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.03"
$ws.Range("E51").Value = "  -1.70%  "

# Rows 44/45: EnergySwap and Decentraland swapped ranking order.
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6459"
$ws.Range("E44").Value = "  +5.30%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.11"
$ws.Range("E45").Value = "  +3.71%  "
